# Updates the cryptocurrency price/volume table with refreshed figures.
# Price-like values in column D are prefixed with a leading apostrophe so
# Excel stores them as literal text (preserving exact formatting such as
# trailing zeros or thousand-separator dots) instead of coercing them to
# numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.851.20"
$ws.Range("E2").Value = "  +1.02%  "
$ws.Range("D3").Value = "1.851.14"
$ws.Range("E3").Value = "  +0.19%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "'335.28"
$ws.Range("E5").Value = "  +0.28%  "
$ws.Range("E6").Value = "  -0.30%  "
$ws.Range("D7").Value = "'0.4657"
$ws.Range("E7").Value = "  +1.06%  "
$ws.Range("D8").Value = "'0.3869"
$ws.Range("E8").Value = "  -0.23%  "
$ws.Range("D9").Value = "'46.83"
$ws.Range("E9").Value = "  +1.82%  "
$ws.Range("E10").Value = "  -0.21%  "
$ws.Range("E11").Value = "  -3.37%  "
$ws.Range("D12").Value = "'21.36"
$ws.Range("E12").Value = "  -0.97%  "
$ws.Range("D13").Value = "1.857.54"
$ws.Range("E13").Value = "  -0.22%  "
$ws.Range("D14").Value = "'5.903"
$ws.Range("E14").Value = "  -1.33%  "
$ws.Range("D15").Value = "'7.161"
$ws.Range("E15").Value = "  -0.02%  "
$ws.Range("D16").Value = "'1.007"
$ws.Range("E16").Value = "  -0.27%  "
$ws.Range("D17").Value = "'90.31"
$ws.Range("E17").Value = "  +2.11%  "
$ws.Range("D18").Value = "'0.06613"
$ws.Range("E18").Value = "  -1.15%  "
$ws.Range("E19").Value = "  -0.76%  "
$ws.Range("D20").Value = "'17.36"
$ws.Range("E20").Value = "  +0.93%  "
$ws.Range("D21").Value = "'1.005"
$ws.Range("E21").Value = "  -0.27%  "
$ws.Range("D22").Value = "27.835.81"
$ws.Range("E22").Value = "  +0.92%  "
$ws.Range("D23").Value = "'5.350"
$ws.Range("E23").Value = "  -1.09%  "
$ws.Range("E24").Value = "  -0.84%  "
$ws.Range("D25").Value = "'2.296"
$ws.Range("E25").Value = "  -0.92%  "
$ws.Range("D26").Value = "2.062.87"
$ws.Range("E26").Value = "  -0.85%  "
$ws.Range("D27").Value = "'158.52"
$ws.Range("E27").Value = "  -0.30%  "
$ws.Range("D28").Value = "'19.48"
$ws.Range("E28").Value = "  -0.33%  "
$ws.Range("D29").Value = "'2.066"
$ws.Range("E29").Value = "  -2.75%  "
$ws.Range("D30").Value = "'5.375"
$ws.Range("E30").Value = "  -1.28%  "
$ws.Range("D31").Value = "'118.86"
$ws.Range("E31").Value = "  -1.49%  "
$ws.Range("D32").Value = "'0.09421"
$ws.Range("E32").Value = "  +0.16%  "
$ws.Range("D33").Value = "'0.9481"
$ws.Range("E33").Value = "  -2.78%  "
$ws.Range("D34").Value = "'3.591"
$ws.Range("E34").Value = "  -0.92%  "
$ws.Range("D35").Value = "'5.266"
$ws.Range("E35").Value = "  -0.69%  "
$ws.Range("D36").Value = "'1.330"
$ws.Range("E36").Value = "  -0.76%  "
$ws.Range("D37").Value = "'0.06026"
$ws.Range("E37").Value = "  +0.08%  "
$ws.Range("D38").Value = "'0.02210"
$ws.Range("E38").Value = "  -0.88%  "
$ws.Range("D39").Value = "'8.267"
$ws.Range("E39").Value = "  -0.49%  "
$ws.Range("D40").Value = "'1.005"
$ws.Range("D41").Value = "'1.159"
$ws.Range("E41").Value = "  -2.29%  "
$ws.Range("D42").Value = "'0.5815"
$ws.Range("E42").Value = "  -1.68%  "
$ws.Range("E43").Value = "  -1.13%  "
$ws.Range("E44").Value = "  -2.79%  "
$ws.Range("D45").Value = "'1.281"
$ws.Range("E45").Value = "  +2.85%  "
$ws.Range("D46").Value = "'0.5455"
$ws.Range("E46").Value = "  -2.55%  "
$ws.Range("D47").Value = "'11.93"
$ws.Range("E47").Value = "  -1.59%  "
$ws.Range("E48").Value = "  +1.01%  "
$ws.Range("D49").Value = "'0.06850"
$ws.Range("E49").Value = "  +2.23%  "
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").Value = "'110.82"
$ws.Range("E50").Value = "  +0.22%  "
$ws.Range("B51").Value = "PaxosStandard"
$ws.Range("C51").Value = "https://coinranking.com/coin/B8xT718SbVhhh+paxosstandard-pax"
$ws.Range("D51").Value = "'1.006"
$ws.Range("E51").Value = "  -32.54%  "
